# Added non-mean centered body awareness: append a row for the new
# "indiv_body.c" variable (its Values/Notes columns are intentionally
# left blank, matching the other derived/mean-centered variable rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "indiv_body.c"
$ws.Range("B40").Value = "indiv_body mean centered"

# Match the formatting used by the other two-column-only rows (e.g. row 36,
# "awareness.c") instead of the default sheet style.
$ws.Range("B36").Copy() | Out-Null
$ws.Range("A40:B40").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection on the newly added row, as in the saved file.
$ws.Range("A40").Select() | Out-Null
